$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old JOIN block (common to all 7 SQL queries on the sheet) that needs
# to be replaced with the corrected join-column block.
$oldBlock = @'
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
'@

$newBlock = @'
LEFT JOIN 
    df_participant prt ON std.study_id = prt."study.study_id"
LEFT JOIN 
    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"
LEFT JOIN 
    df_treatments trt ON prt.participant_id = trt."participant.participant_id"
LEFT JOIN 
    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"
LEFT JOIN 
    df_survival srv ON prt.participant_id = srv."participant.participant_id"
LEFT JOIN 
    df_reference_files rfs ON std.study_id = rfs."study.study_id"
'@

# Every query cell that contains the old join block.
$queryCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $queryCells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value2
    if ($text -ne $null -and $text.Contains($oldBlock)) {
        $rng.Value = $text.Replace($oldBlock, $newBlock)
    }
}

# Resize column C to a fixed width of 68 (no longer "best fit").
$ws.Columns.Item(3).ColumnWidth = 67.1666666666667

# Update the selected / active cell from C7 to B2 and scroll the view back
# to the top-left (A1) instead of A6.
$ws.Range("A1").Select() | Out-Null
$ws.Range("B2").Select() | Out-Null
